$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G across rows 2-7
# (Column F is left unchanged; G = B + C + D + E)
$data = @{
    2 = @{ B = 0.7287194209349384;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 3.594575437922795  }
    3 = @{ B = 3.182878228561681;   C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729  }
    4 = @{ B = 3.182878228561681;   C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538  }
    5 = @{ B = 3.182878228561681;   C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538  }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538  }
    7 = @{ B = 0.7287194209349384;  C = 0.3375848360084654; D = 16.98373111632243;   E = 0.4998867070740569; G = 18.54992208033989  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
